$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9: sample[1]
$ws.Range("A9").Value = "sample[1]"

# Copy A8's formatting (font + alignment) onto A9 to match the other label cells
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B9").Value = 3
$ws.Range("C9").Formula = "=INT(32767*SIN(2*PI()*B9/(`$B`$3/`$B`$1)))"

# Update the selection to mimic the saved workbook state
$ws.Range("B10").Select()
